$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '71.954.72'
Set-TextValue 'E2' '  +0.14%  '
Set-TextValue 'D3' '2.685.10'
Set-TextValue 'E3' '  +2.07%  '
Set-TextValue 'E4' '  +0.08%  '
Set-TextValue 'D5' '597.52'
Set-TextValue 'E5' '  -1.83%  '
Set-TextValue 'E6' '  -2.90%  '
Set-TextValue 'E7' '  +0.06%  '
Set-TextValue 'D8' '0.524'
Set-TextValue 'E8' '  -0.08%  '
Set-TextValue 'D9' '2.683.51'
Set-TextValue 'E9' '  +2.13%  '
Set-TextValue 'E10' '  -1.80%  '
Set-TextValue 'E11' '  +2.13%  '
Set-TextValue 'D12' '0.354'
Set-TextValue 'E12' '  +1.71%  '
Set-TextValue 'D13' '4.99'
Set-TextValue 'E13' '  -0.66%  '
Set-TextValue 'D14' '3.176.85'
Set-TextValue 'E14' '  +3.00%  '
Set-TextValue 'D15' '71.842.85'
Set-TextValue 'E15' '  +0.10%  '
Set-TextValue 'E16' '  -2.07%  '
Set-TextValue 'D17' '26.18'
Set-TextValue 'E17' '  -1.26%  '
Set-TextValue 'D18' '2.682.09'
Set-TextValue 'E18' '  +2.05%  '
Set-TextValue 'D19' '12.20'
Set-TextValue 'E19' '  +6.10%  '
Set-TextValue 'D20' '8.13'
Set-TextValue 'E20' '  +1.18%  '
Set-TextValue 'D21' '371.42'
Set-TextValue 'E21' '  -3.00%  '
Set-TextValue 'D22' '4.19'
Set-TextValue 'E22' '  +0.92%  '
Set-TextValue 'E23' '  -0.53%  '
Set-TextValue 'D24' '72.24'
Set-TextValue 'E24' '  -0.62%  '
Set-TextValue 'E25' '  +0.02%  '
Set-TextValue 'E26' '  -2.23%  '
Set-TextValue 'D27' '9.78'
Set-TextValue 'E27' '  -1.55%  '
Set-TextValue 'D28' '2.821.45'
Set-TextValue 'E28' '  +2.14%  '
Set-TextValue 'D29' '0.996'
Set-TextValue 'E29' '  -0.68%  '
Set-TextValue 'D30' '0.0₃0961'
Set-TextValue 'E30' '  -0.55%  '
Set-TextValue 'E31' '  +0.03%  '
Set-TextValue 'D32' '498.57'
Set-TextValue 'E32' '  -8.75%  '
Set-TextValue 'E33' '  -2.32%  '
Set-TextValue 'D34' '1.82'
Set-TextValue 'E34' '  -0.77%  '
Set-TextValue 'E35' '  +0.06%  '
Set-TextValue 'D36' '163.91'
Set-TextValue 'E36' '  -1.34%  '
Set-TextValue 'D37' '19.58'
Set-TextValue 'E37' '  +1.88%  '
Set-TextValue 'D38' '19.10'
Set-TextValue 'E38' '  -0.11%  '
Set-TextValue 'E39' '  -1.24%  '
Set-TextValue 'E40' '  -6.36%  '
Set-TextValue 'D41' '1.78'
Set-TextValue 'E41' '  -3.92%  '
Set-TextValue 'E42' '  -0.01%  '
Set-TextValue 'D43' '5.01'
Set-TextValue 'E43' '  -0.41%  '
Set-TextValue 'D44' '0.334'
Set-TextValue 'E44' '  +0.11%  '
Set-TextValue 'D45' '2.55'
Set-TextValue 'E45' '  -2.34%  '
Set-TextValue 'D46' '157.69'
Set-TextValue 'E46' '  +4.68%  '
Set-TextValue 'D47' '39.39'
Set-TextValue 'E47' '  -0.16%  '
Set-TextValue 'D48' '0.565'
Set-TextValue 'E48' '  +5.03%  '
Set-TextValue 'D49' '3.74'
Set-TextValue 'E49' '  +1.97%  '
Set-TextValue 'D50' '1.76'
Set-TextValue 'E50' '  +3.98%  '
Set-TextValue 'D51' '0.0762'
Set-TextValue 'E51' '  +0.85%  '
